# Auto-generated edit script: updates 2025 (column L) violent-crime counts
# for "2025-06-07" data refresh across Citywide Totals, By Neighborhood, and
# each affected neighborhood sheet, per the commit diff.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2731
$ws.Range("L3").Value = 2777
$ws.Range("L4").Value = 739
$ws.Range("L6").Value = 2478
$ws.Range("L7").Value = 8885

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 293
$ws.Range("L8").Value = 560
$ws.Range("L9").Value = 57
$ws.Range("L11").Value = 155
$ws.Range("L18").Value = 62
$ws.Range("L20").Value = 227
$ws.Range("L23").Value = 91
$ws.Range("L29").Value = 478
$ws.Range("L31").Value = 84
$ws.Range("L33").Value = 401
$ws.Range("L34").Value = 53
$ws.Range("L37").Value = 327
$ws.Range("L40").Value = 24
$ws.Range("L42").Value = 294
$ws.Range("L44").Value = 66
$ws.Range("L47").Value = 67
$ws.Range("L49").Value = 48
$ws.Range("L50").Value = 47
$ws.Range("L51").Value = 108
$ws.Range("L52").Value = 179
$ws.Range("L55").Value = 86
$ws.Range("L60").Value = 54
$ws.Range("L63").Value = 27
$ws.Range("L65").Value = 163
$ws.Range("L66").Value = 21
$ws.Range("L67").Value = 328
$ws.Range("L68").Value = 28
$ws.Range("L72").Value = 39
$ws.Range("L73").Value = 74
$ws.Range("L76").Value = 109
$ws.Range("L78").Value = 116
$ws.Range("L79").Value = 238
$ws.Range("L83").Value = 210
$ws.Range("L85").Value = 459
$ws.Range("L90").Value = 89
$ws.Range("L91").Value = 127
$ws.Range("L92").Value = 26
$ws.Range("L94").Value = 108
$ws.Range("L95").Value = 121
$ws.Range("L96").Value = 87
$ws.Range("L99").Value = 148
$ws.Range("L101").Value = 8885

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 87

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 89
$ws.Range("L7").Value = 293

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 56
$ws.Range("L7").Value = 155

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 187
$ws.Range("L7").Value = 459

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 179

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 186
$ws.Range("L7").Value = 560

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 66
$ws.Range("L7").Value = 210

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 126
$ws.Range("L6").Value = 138
$ws.Range("L7").Value = 401

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 121

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 327

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 163

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 36
$ws.Range("L7").Value = 148

# Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 84

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 119
$ws.Range("L7").Value = 328

# Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 48

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 153
$ws.Range("L3").Value = 178
$ws.Range("L6").Value = 121
$ws.Range("L7").Value = 478

# Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 66

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 109

# Ashburn
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 29
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 70

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 84
$ws.Range("L3").Value = 89
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 294

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 116

# Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 86

# Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 28
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 91

# Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 127

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 86
$ws.Range("L7").Value = 238

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 66
$ws.Range("L7").Value = 227

# Calumet Heights
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 62

# Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 53

# West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 108

# Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 67

# Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L4").Value = 5
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 47

# North Center
$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 21

# Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 57

# Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 32
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 74

# West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 26

# Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 26
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 89

# Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 30
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 108

# North Park
$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 28

# Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 54

# Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 39

# Hegewisch
$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 24

